$wb = $excel.ActiveWorkbook

# New row (row 68) to append at the bottom of each of the 4 sheets.
$newRows = @(
    @{ Sheet = "ROW35-FE-LIFTER";  A = "2025-03-07 03:42:06"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"; E = "0x d"; G = "568631262647113770877196"; I = 13 },
    @{ Sheet = "ROW35-MID-LIFTER"; A = "2025-03-07 03:29:35"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; E = "0x e"; G = "568631262647113770942732"; I = 14 },
    @{ Sheet = "ROW02-FE-LIFTER";  A = "2025-03-07 03:51:45"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"; E = "0xff"; G = "568631262647113769959692"; I = 255 },
    @{ Sheet = "ROW02-MID-LIFTER"; A = "2025-03-07 03:41:15"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; E = "0x 3"; G = "568631262647113769959692"; I = 3 }
)

foreach ($item in $newRows) {
    $ws = $wb.Worksheets.Item($item.Sheet)
    $r = 68

    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = "0x01,0x90 "
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = "0x01,0x90,"
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = 400

    # Column G holds a 24-digit numeric string. Force it to be stored as
    # text (like the rest of the column) instead of being coerced into a
    # floating point number with loss of precision, then restore the
    # default "Normal" style so no extraneous formatting is left behind.
    $gCell = $ws.Cells.Item($r, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $item.G
    $gCell.Style = "Normal"

    $ws.Cells.Item($r, 8).Value = 400
    $ws.Cells.Item($r, 9).Value = $item.I
}
